# Auto-generated edit script applying cached-value updates to the
# Leve profit calculator sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR),
# per the scheduled market-data refresh described in the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2961.6155
$ws.Range("J51").Value = 2870.1
$ws.Range("L51").Value = 2870.1
$ws.Range("N51").Value = -3838.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 772.95
$ws.Range("I53").Value = 353.18182
$ws.Range("K53").Value = 353.18182
$ws.Range("M53").Value = 283.81818

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 19624
$ws.Range("J95").Value = 19624
$ws.Range("L95").Value = 19624
$ws.Range("N95").Value = -25116

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 839.6316
$ws.Range("I127").Value = 445.7
$ws.Range("J127").Value = 1277.3334
$ws.Range("K127").Value = 1337.1
$ws.Range("L127").Value = 3832.0002
$ws.Range("M127").Value = 3622.9
$ws.Range("N127").Value = -13752.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 194560.39
$ws.Range("I132").Value = 215151.88
$ws.Range("J132").Value = 1000.4
$ws.Range("K132").Value = 645455.64
$ws.Range("L132").Value = 3001.2
$ws.Range("M132").Value = -642925.64
$ws.Range("N132").Value = -8061.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27028720
$ws.Range("I137").Value = 1164.3793
$ws.Range("J137").Value = 125003620
$ws.Range("K137").Value = 3493.1379
$ws.Range("L137").Value = 375010860
$ws.Range("M137").Value = -943.1379000000002
$ws.Range("N137").Value = -375015960

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2177146
$ws.Range("I138").Value = 3848722
$ws.Range("J138").Value = 4097.25
$ws.Range("K138").Value = 11546166
$ws.Range("L138").Value = 12291.75
$ws.Range("M138").Value = -11541026
$ws.Range("N138").Value = -22571.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1373
$ws.Range("I61").Value = 1397.1052
$ws.Range("J61").Value = 1296.6666
$ws.Range("K61").Value = 1397.1052
$ws.Range("L61").Value = 1296.6666
$ws.Range("M61").Value = -1185.1052
$ws.Range("N61").Value = -1720.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 28400
$ws.Range("J109").Value = 28400
$ws.Range("L109").Value = 28400
$ws.Range("N109").Value = -31174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1373
$ws.Range("I136").Value = 1397.1052
$ws.Range("J136").Value = 1296.6666
$ws.Range("K136").Value = 4191.3156
$ws.Range("L136").Value = 3889.9998
$ws.Range("M136").Value = -1641.3156
$ws.Range("N136").Value = -8989.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15834.637
$ws.Range("I82").Value = 7282.8
$ws.Range("J82").Value = 22961.166
$ws.Range("K82").Value = 7282.8
$ws.Range("L82").Value = 22961.166
$ws.Range("M82").Value = -6899.8
$ws.Range("N82").Value = -23727.166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 15834.637
$ws.Range("I85").Value = 7282.8
$ws.Range("J85").Value = 22961.166
$ws.Range("K85").Value = 7282.8
$ws.Range("L85").Value = 22961.166
$ws.Range("M85").Value = -5956.8
$ws.Range("N85").Value = -25613.166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1294.2
$ws.Range("I107").Value = 1262.9231
$ws.Range("K107").Value = 1262.9231
$ws.Range("M107").Value = 657.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1647.2858
$ws.Range("I31").Value = 1374.0952
$ws.Range("J31").Value = 2057.0715
$ws.Range("K31").Value = 1374.0952
$ws.Range("L31").Value = 2057.0715
$ws.Range("M31").Value = -1079.0952
$ws.Range("N31").Value = -2647.0715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1647.2858
$ws.Range("I34").Value = 1374.0952
$ws.Range("J34").Value = 2057.0715
$ws.Range("K34").Value = 1374.0952
$ws.Range("L34").Value = 2057.0715
$ws.Range("M34").Value = -1172.0952
$ws.Range("N34").Value = -2461.0715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9800.182000000001
$ws.Range("J51").Value = 9870.200000000001
$ws.Range("L51").Value = 9870.200000000001
$ws.Range("N51").Value = -11342.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 9363.25
$ws.Range("J60").Value = 10484.333
$ws.Range("L60").Value = 10484.333
$ws.Range("N60").Value = -11506.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9800.182000000001
$ws.Range("J61").Value = 9870.200000000001
$ws.Range("L61").Value = 9870.200000000001
$ws.Range("N61").Value = -10566.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 22147.5
$ws.Range("J68").Value = 22147.5
$ws.Range("L68").Value = 22147.5
$ws.Range("N68").Value = -23645.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 22147.5
$ws.Range("J71").Value = 22147.5
$ws.Range("L71").Value = 66442.5
$ws.Range("N71").Value = -73930.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 17382
$ws.Range("J74").Value = 17382
$ws.Range("L74").Value = 17382
$ws.Range("N74").Value = -19130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 17382
$ws.Range("J77").Value = 17382
$ws.Range("L77").Value = 52146
$ws.Range("N77").Value = -60882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 7703.6665
$ws.Range("I107").Value = 7703.6665
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 7703.6665
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -5783.6665
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 59702
$ws.Range("J111").Value = 59702
$ws.Range("L111").Value = 59702
$ws.Range("N111").Value = -67882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -3338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -3204

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1162.1666
$ws.Range("I44").Value = 435.14285
$ws.Range("J44").Value = 2180
$ws.Range("K44").Value = 1305.42855
$ws.Range("L44").Value = 6540
$ws.Range("M44").Value = -907.4285500000001
$ws.Range("N44").Value = -7336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3611.9722
$ws.Range("I131").Value = 5921.222
$ws.Range("J131").Value = 2842.2222
$ws.Range("K131").Value = 17763.666
$ws.Range("L131").Value = 8526.6666
$ws.Range("M131").Value = -12723.666
$ws.Range("N131").Value = -18606.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 397.22223
$ws.Range("I107").Value = 255.3125
$ws.Range("J107").Value = 603.63635
$ws.Range("K107").Value = 255.3125
$ws.Range("L107").Value = 603.63635
$ws.Range("M107").Value = 1664.6875
$ws.Range("N107").Value = -4443.63635

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 64753.332
$ws.Range("J140").Value = 64753.332
$ws.Range("L140").Value = 64753.332
$ws.Range("N140").Value = -75113.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1753.4857
$ws.Range("I68").Value = 1670.4286
$ws.Range("J68").Value = 2085.7144
$ws.Range("K68").Value = 1670.4286
$ws.Range("L68").Value = 2085.7144
$ws.Range("M68").Value = -921.4286
$ws.Range("N68").Value = -3583.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1753.4857
$ws.Range("I71").Value = 1670.4286
$ws.Range("J71").Value = 2085.7144
$ws.Range("K71").Value = 8352.143
$ws.Range("L71").Value = 10428.572
$ws.Range("M71").Value = -4608.143
$ws.Range("N71").Value = -17916.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2348.8845
$ws.Range("I136").Value = 1428.55
$ws.Range("J136").Value = 5416.6665
$ws.Range("K136").Value = 4285.65
$ws.Range("L136").Value = 16249.9995
$ws.Range("M136").Value = -1735.65
$ws.Range("N136").Value = -21349.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 604.5454999999999
$ws.Range("I107").Value = 231.25
$ws.Range("K107").Value = 693.75
$ws.Range("M107").Value = 1226.25
